# Updated 24V test data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels and Devices")

# Update the expected 24V PSU load values in row 8
$ws.Range("F8").Value = 0.205
$ws.Range("J8").Value = 0.207
$ws.Range("N8").Value = "'0.208"
$ws.Range("O8").Value = 0.207

# Update the view/selection state to match the last-saved selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("O10").Select()
